$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Added skybox: row 23 (Infinite Sky Box) marked complete on Milestone II ---
$ws.Range("E23").Value = "II"
$ws.Range("F23").Value = "X"

# --- Milestone re-assignment: row 33 (Combining 2 functional lights) moved from Milestone II to Milestone I ---
$ws.Range("E33").Value = "I"

# --- Added multithreading: rows 55 & 56 marked complete on Milestone II ---
$ws.Range("E55").Value = "II"
$ws.Range("F55").Value = "X"

$ws.Range("E56").Value = "II"
$ws.Range("F56").Value = "X"

# --- New project source citation for model loading tutorial ---
$ws.Range("A95").Value = "http://www.opengl-tutorial.org/beginners-tutorials/tutorial-7-model-loading/"

# --- Scroll the view down so row 52 is at the top (selection remains E33) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1

$wb.Application.CalculateFull()
